# Refresh the scraped "cryptos" price/volume snapshot (Sun Mar 19 23:42:53
# UTC 2023 GitHub Actions run). Columns: B=Coin, C=Link, D=Price,
# E=Volume(1h) for data rows 2..51. A handful of rows also swap which coin
# occupies that rank (26/28 ImmutableX<->EthereumClassic, 38/39
# TheSandbox<->Hedera), so B/C get rewritten there too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.244.77"
$ws.Cells.Item(2, 5).Value = "  +4.12%  "

$ws.Cells.Item(3, 4).Value = "1.809.05"
$ws.Cells.Item(3, 5).Value = "  +2.04%  "

$ws.Cells.Item(4, 4).Value = "1.004"
$ws.Cells.Item(4, 5).Value = "  -0.52%  "

$ws.Cells.Item(5, 4).Value = "338.99"
$ws.Cells.Item(5, 5).Value = "  +1.65%  "

$ws.Cells.Item(6, 4).Value = "'1.000"
$ws.Cells.Item(6, 5).Value = "  -0.33%  "

$ws.Cells.Item(7, 4).Value = "0.3918"
$ws.Cells.Item(7, 5).Value = "  +4.07%  "

$ws.Cells.Item(8, 4).Value = "0.3502"
$ws.Cells.Item(8, 5).Value = "  +3.02%  "

$ws.Cells.Item(9, 4).Value = "48.59"
$ws.Cells.Item(9, 5).Value = "  +1.08%  "

$ws.Cells.Item(10, 4).Value = "1.182"

$ws.Cells.Item(11, 4).Value = "0.07567"
$ws.Cells.Item(11, 5).Value = "  +2.63%  "

$ws.Cells.Item(12, 4).Value = "1.001"
$ws.Cells.Item(12, 5).Value = "  -0.53%  "

$ws.Cells.Item(13, 4).Value = "22.17"
$ws.Cells.Item(13, 5).Value = "  +3.84%  "

$ws.Cells.Item(14, 4).Value = "6.543"
$ws.Cells.Item(14, 5).Value = "  +2.87%  "

$ws.Cells.Item(15, 4).Value = "1.813.60"
$ws.Cells.Item(15, 5).Value = "  +2.04%  "

$ws.Cells.Item(16, 4).Value = "7.182"
$ws.Cells.Item(16, 5).Value = "  +2.91%  "

$ws.Cells.Item(17, 4).Value = "0.00001108"
$ws.Cells.Item(17, 5).Value = "  +2.63%  "

$ws.Cells.Item(18, 4).Value = "0.06718"
$ws.Cells.Item(18, 5).Value = "  +1.12%  "

$ws.Cells.Item(19, 4).Value = "85.38"
$ws.Cells.Item(19, 5).Value = "  +2.18%  "

$ws.Cells.Item(20, 4).Value = "0.9996"

$ws.Cells.Item(21, 4).Value = "17.79"
$ws.Cells.Item(21, 5).Value = "  +4.06%  "

$ws.Cells.Item(22, 4).Value = "6.586"
$ws.Cells.Item(22, 5).Value = "  +1.15%  "

$ws.Cells.Item(23, 4).Value = "28.253.17"
$ws.Cells.Item(23, 5).Value = "  +4.07%  "

$ws.Cells.Item(24, 4).Value = "12.49"
$ws.Cells.Item(24, 5).Value = "  +1.48%  "

$ws.Cells.Item(25, 4).Value = "2.405"
$ws.Cells.Item(25, 5).Value = "  -0.81%  "

$ws.Cells.Item(26, 2).Value = "EthereumClassic"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(26, 4).Value = "21.46"
$ws.Cells.Item(26, 5).Value = "  +2.52%  "

$ws.Cells.Item(27, 5).Value = "  +1.61%  "

$ws.Cells.Item(28, 2).Value = "ImmutableX"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(28, 4).Value = "'1.480"
$ws.Cells.Item(28, 5).Value = "  -0.50%  "

$ws.Cells.Item(29, 4).Value = "154.94"
$ws.Cells.Item(29, 5).Value = "  +2.48%  "

$ws.Cells.Item(30, 4).Value = "2.019.91"
$ws.Cells.Item(30, 5).Value = "  +2.07%  "

$ws.Cells.Item(31, 4).Value = "136.77"
$ws.Cells.Item(31, 5).Value = "  +3.35%  "

$ws.Cells.Item(32, 4).Value = "6.435"
$ws.Cells.Item(32, 5).Value = "  +8.92%  "

$ws.Cells.Item(33, 4).Value = "4.024"
$ws.Cells.Item(33, 5).Value = "  -0.83%  "

$ws.Cells.Item(34, 4).Value = "'0.08860"
$ws.Cells.Item(34, 5).Value = "  +2.92%  "

$ws.Cells.Item(35, 4).Value = "13.22"
$ws.Cells.Item(35, 5).Value = "  +2.76%  "

$ws.Cells.Item(36, 4).Value = "0.02469"
$ws.Cells.Item(36, 5).Value = "  +6.50%  "

$ws.Cells.Item(37, 4).Value = "5.502"
$ws.Cells.Item(37, 5).Value = "  +2.73%  "

$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(38, 4).Value = "'0.06560"
$ws.Cells.Item(38, 5).Value = "  +4.49%  "

$ws.Cells.Item(39, 2).Value = "TheSandbox"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(39, 4).Value = "0.6936"
$ws.Cells.Item(39, 5).Value = "  +2.76%  "

$ws.Cells.Item(40, 4).Value = "1.619"
$ws.Cells.Item(40, 5).Value = "  -1.84%  "

$ws.Cells.Item(41, 4).Value = "0.2233"
$ws.Cells.Item(41, 5).Value = "  +3.13%  "

$ws.Cells.Item(42, 5).Value = "  +3.08%  "

$ws.Cells.Item(43, 4).Value = "8.558"
$ws.Cells.Item(43, 5).Value = "  -1.21%  "

$ws.Cells.Item(44, 4).Value = "14.69"
$ws.Cells.Item(44, 5).Value = "  +2.53%  "

$ws.Cells.Item(45, 5).Value = "  +2.62%  "

$ws.Cells.Item(46, 4).Value = "0.9994"
$ws.Cells.Item(46, 5).Value = "  -0.46%  "

$ws.Cells.Item(47, 4).Value = "3.876"
$ws.Cells.Item(47, 5).Value = "  +1.31%  "

$ws.Cells.Item(48, 4).Value = "2.171"
$ws.Cells.Item(48, 5).Value = "  +3.51%  "

$ws.Cells.Item(49, 4).Value = "132.02"
$ws.Cells.Item(49, 5).Value = "  +2.80%  "

$ws.Cells.Item(50, 4).Value = "'0.07250"
$ws.Cells.Item(50, 5).Value = "  +1.69%  "

$ws.Cells.Item(51, 4).Value = "80.55"
$ws.Cells.Item(51, 5).Value = "  +2.87%  "
